$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 264
$ws.Range("F4").Value = 28
$ws.Range("F8").Value = 163
$ws.Range("F10").Value = 230
$ws.Range("F11").Value = 6086
$ws.Range("F21").Value = 713
$ws.Range("F22").Value = 160
$ws.Range("F25").Value = 1025
$ws.Range("F27").Value = 1842
$ws.Range("F28").Value = 502
$ws.Range("F29").Value = 33

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 261

# Sheet "全部类型" (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 261
$ws.Range("F4").Value = 264
$ws.Range("F5").Value = 28
$ws.Range("F10").Value = 163
$ws.Range("F12").Value = 230
$ws.Range("F13").Value = 6086
$ws.Range("F28").Value = 713
$ws.Range("F32").Value = 160
$ws.Range("F35").Value = 1025
$ws.Range("F37").Value = 1842
$ws.Range("F38").Value = 502
$ws.Range("F39").Value = 33

$wb.Save()
